# Updated cryptos list values (prices + 1h volume %) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.848.94'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '3.742.33'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.95'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.22'
$ws.Range("E6").Value = '  +1.70%  '
$ws.Range("D7").Value = '3.740.69'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  -2.33%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.60'
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("E12").Value = '  -4.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.97'
$ws.Range("E13").Value = '  -2.51%  '
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = '4.366.00'
$ws.Range("D16").Value = '3.745.66'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '69.920.04'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.55'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.53'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '503.71'
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.19'
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.724'
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.62'
$ws.Range("E24").Value = '  +5.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.52'
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.35'
$ws.Range("E26").Value = '  +4.13%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '13.02'
$ws.Range("E27").Value = '  -3.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000135'
$ws.Range("E28").Value = '  +7.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.92'
$ws.Range("E31").Value = '  +2.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.99'
$ws.Range("E32").Value = '  +1.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.67'
$ws.Range("E33").Value = '  -2.93%  '
$ws.Range("E34").Value = '  -1.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("E38").Value = '  +2.62%  '
$ws.Range("E39").Value = '  +5.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.08'
$ws.Range("E40").Value = '  +12.84%  '
$ws.Range("E41").Value = '  -5.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.87'
$ws.Range("E42").Value = '  +2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '438.90'
$ws.Range("E43").Value = '  +4.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '49.66'
$ws.Range("E44").Value = '  -3.46%  '
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("D46").Value = '2.973.62'
$ws.Range("E46").Value = '  -3.36%  '
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.31'
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.29'
$ws.Range("E49").Value = '  +1.96%  '
$ws.Range("E51").Value = '  -1.55%  '
